$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated TPM-derived values
$ws.Range("M2").Value = 0.295943
$ws.Range("N2").Value = 0.887829
$ws.Range("O2").Value = 0.0553754985456454
$ws.Range("P2").Value = 0.0553754985456454
$ws.Range("Q2").Value = 0.005257131452
$ws.Range("R2").Value = 0.047314183068
$ws.Range("S2").Value = 0.0553754985456454
$ws.Range("T2").Value = 0.0553754985456454

# Row 3 updated TPM-derived values
$ws.Range("O3").Value = 0.9446245014543545
$ws.Range("P3").Value = 0.9446245014543546
$ws.Range("S3").Value = 0.9446245014543545
$ws.Range("T3").Value = 0.9446245014543546
